# Actualización automática hashcode mié jul  3 01:30:27 CEST 2019
# Updates the "hashcode" column (B) for a set of rows in the active sheet,
# matching the shared-string replacements from the canonical OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B100").Value = "85819c9b0ee864700a6fb3abf7b62758"
$ws.Range("B104").Value = "afc45b0ea45fcd2114d8102997488408"
$ws.Range("B113").Value = "956b266fd844e9f3fca2194ee278fadb"
$ws.Range("B122").Value = "d15ca3c8fb72fbbd22db7c2394f28a69"
$ws.Range("B164").Value = "0a80cf60deec27272e68c8141fbee685"
$ws.Range("B175").Value = "c2c7e5ddfc176b5bfb3788ec418d8ad3"
$ws.Range("B227").Value = "b025351c83f75e31b423737f214610be"
$ws.Range("B230").Value = "a7ccd9496d18261177551264266f67e7"
$ws.Range("B232").Value = "1547e149e4f8ab0c718a75fa7757d8df"
$ws.Range("B233").Value = "380c5e4c6ed05e85df43317f9a0cfa66"
$ws.Range("B331").Value = "d9986ed4380897b50d61c0803314de7c"
$ws.Range("B342").Value = "052d5b4453144717d9154004c40aed09"
$ws.Range("B343").Value = "9c8e173b79f48d63f00af95644862e76"
$ws.Range("B381").Value = "426758b07b194188b97fe09b886f440d"
$ws.Range("B458").Value = "62f05aaa5756711c583f9c74bdffd409"
$ws.Range("B477").Value = "e1b8840a7130774ea1c4a2335241f85b"
$ws.Range("B478").Value = "c43271c014c176323131e06d59c1a61b"
$ws.Range("B619").Value = "bd09cfb4e9f5a5a1edc58ee2f6cbef23"
$ws.Range("B623").Value = "5df9e1ffb7ca51b90d6720532ccfee6f"
$ws.Range("B628").Value = "ae8a27b09551a4de674da30e82a0e23c"
$ws.Range("B649").Value = "fb00f9d8dcf9078f088babf8d5927fb5"
$ws.Range("B655").Value = "5d69577e6a8352311ecdf23d993300dd"
$ws.Range("B779").Value = "babf3fd530aff2ea45435a4292853ff1"
$ws.Range("B818").Value = "4c2ed9e49577e877cba8646fab52dc00"
$ws.Range("B831").Value = "3ebef27ff7385eb5bb0c6c1d9dc07834"
$ws.Range("B874").Value = "c9c849f03081bb7a17b5eba5feebb7ea"
